$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell assignments. Each text value is written with a
# leading apostrophe to force literal text (avoids Excel auto-converting
# numeric-looking strings such as "595.06" or "1.00" into real numbers),
# then ClearFormats() strips the resulting "text" number-format / quote-prefix
# style so the cell keeps its original (unstyled) appearance.

$ws.Range("D2").Value = "'67.898.20"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.47%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.618.95"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -0.72%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'595.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.68%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'152.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -1.26%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.544"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -0.54%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'2.618.31"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -0.73%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.133"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +6.49%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -0.70%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'5.19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -0.83%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.348"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -1.25%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'27.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -1.78%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.0000188"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +1.80%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.097.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -0.50%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'67.698.93"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +0.20%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.613.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -0.97%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'372.29"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +1.54%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'11.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.56%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'7.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -2.66%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'4.22"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -2.09%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'4.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -3.25%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'2.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -4.30%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'72.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +9.57%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E27").Value = "'  -1.66%  "
$ws.Range("E27").ClearFormats()
$ws.Range("B28").Value = "'WrappedeETH"
$ws.Range("B28").ClearFormats()
$ws.Range("C28").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("C28").ClearFormats()
$ws.Range("D28").Value = "'2.757.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +0.06%  "
$ws.Range("E28").ClearFormats()
$ws.Range("B29").Value = "'PEPE"
$ws.Range("B29").ClearFormats()
$ws.Range("C29").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C29").ClearFormats()
$ws.Range("D29").Value = "'0.0000103"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -1.16%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'588.22"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +0.77%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -0.33%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'7.78"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -2.05%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'1.37"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -3.52%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'1.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -1.28%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +0.05%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.126"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -3.03%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'1.52"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -1.24%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'157.77"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -0.26%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'19.10"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -2.00%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'1.89"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +2.58%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.367"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -1.43%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'5.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -1.38%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'2.68"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +2.17%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'17.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +4.62%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'  +0.03%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'40.42"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -1.79%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'155.98"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.21%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.0₆0295"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +1.80%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'3.67"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -2.02%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'1.69"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -2.81%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.0778"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.23%  "
$ws.Range("E51").ClearFormats()
